# Weekly update: two new price records were collected and inserted at the
# top of the data block (rows 464-465), pushing the existing records
# (previously rows 464-571) down by two rows (now rows 466-573).
#
# We copy each row's A:R values from its old position (r-2) into its new
# position (r), working from the bottom of the sheet upward so that we
# never overwrite a source row before it has been read.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 573; $r -ge 466; $r--) {
    $src = $r - 2
    $srcRange = $ws.Range("A" + $src + ":R" + $src)
    $dstRange = $ws.Range("A" + $r + ":R" + $r)
    $dstRange.Value2 = $srcRange.Value2
}

# Make sure the date column on the two newly created rows keeps the same
# date/time number format used throughout column D.
$ws.Cells.Item(572, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(573, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 464: brand-new record.
$ws.Cells.Item(464, 1).Value2 = 10
$ws.Cells.Item(464, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(464, 3).Value2 = "La Araucanía"
$ws.Cells.Item(464, 4).Value2 = 45211
$ws.Cells.Item(464, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(464, 5).Value2 = 9
$ws.Cells.Item(464, 6).Value2 = 100112017
$ws.Cells.Item(464, 7).Value2 = "Apio"
$ws.Cells.Item(464, 8).Value2 = "Americana (o)"
$ws.Cells.Item(464, 9).Value2 = "Primera"
$ws.Cells.Item(464, 10).Value2 = 300
$ws.Cells.Item(464, 11).Value2 = 8000
$ws.Cells.Item(464, 12).Value2 = 8000
$ws.Cells.Item(464, 13).Value2 = 8000
$ws.Cells.Item(464, 14).Value2 = "`$/caja 8 unidades"
$ws.Cells.Item(464, 15).Value2 = "Provincia del Elquí"
$ws.Cells.Item(464, 16).Value2 = 8000
$ws.Cells.Item(464, 17).Value2 = 1
$ws.Cells.Item(464, 18).Value2 = "Hortaliza"

# Row 465: brand-new record.
$ws.Cells.Item(465, 1).Value2 = 10
$ws.Cells.Item(465, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(465, 3).Value2 = "La Araucanía"
$ws.Cells.Item(465, 4).Value2 = 45211
$ws.Cells.Item(465, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(465, 5).Value2 = 9
$ws.Cells.Item(465, 6).Value2 = 100112017
$ws.Cells.Item(465, 7).Value2 = "Apio"
$ws.Cells.Item(465, 8).Value2 = "Americana (o)"
$ws.Cells.Item(465, 9).Value2 = "Primera"
$ws.Cells.Item(465, 10).Value2 = 180
$ws.Cells.Item(465, 11).Value2 = 8000
$ws.Cells.Item(465, 12).Value2 = 9000
$ws.Cells.Item(465, 13).Value2 = 8556
$ws.Cells.Item(465, 14).Value2 = "`$/docena de matas"
$ws.Cells.Item(465, 15).Value2 = "Provincia del Elquí"
$ws.Cells.Item(465, 16).Value2 = 1426
$ws.Cells.Item(465, 17).Value2 = 6
$ws.Cells.Item(465, 18).Value2 = "Hortaliza"

Write-Host "Edit complete"
